$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushing existing rows 3..11 down to 4..12)
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the missing "9:00 - 10:00" time slot and the
# same green style used by the other time-slot rows.
$ws.Range("A3").Value = "9:00 - 10:00"

$ws.Range("B3:H3").Value = $null
$ws.Range("B4:H4").Copy()
$ws.Range("B3:H3").PasteSpecial(-4122)

# Update the selection to match the post-edit state.
$ws.Range("C20").Select()
